$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = 'Datos actualizados a 6 de Agosto de 2020 a las 22:35'

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 5016633
$ws.Range("C4").Value = 43065
$ws.Range("D4").Value = 2566837
$ws.Range("E4").Value = 2287391
$ws.Range("G4").Value = 804
$ws.Range("H4").Value = 162405

# Row 22: Alemania -> Alemania
$ws.Range("B22").Value = 215210
$ws.Range("C22").Value = 1106
$ws.Range("E22").Value = 9758

# Row 36: Israel -> Israel
$ws.Range("B36").Value = 79559
$ws.Range("C36").Value = 1640
$ws.Range("D36").Value = 53427
$ws.Range("E36").Value = 25556

# Row 70: Etiopia -> Costa Rica
$ws.Range("A70").Value = 'Costa Rica'
$ws.Range("B70").Value = 21070
$ws.Range("C70").Value = 653
$ws.Range("D70").Value = 7038
$ws.Range("E70").Value = 13832
$ws.Range("H70").Value = 200

# Row 71: Costa Rica -> Etiopia
$ws.Range("A71").Value = 'Etiopia'
$ws.Range("B71").Value = 20900
$ws.Range("C71").Value = 564
$ws.Range("D71").Value = 9027
$ws.Range("E71").Value = 11508
$ws.Range("G71").Value = 9
$ws.Range("H71").Value = 365

# Row 76: Costa de Marfil -> Costa de Marfil
$ws.Range("B76").Value = 16447
$ws.Range("C76").Value = 98
$ws.Range("D76").Value = 12484
$ws.Range("E76").Value = 3860

# Row 89: Guayana Francesa -> Guayana Francesa
$ws.Range("B89").Value = 8127
$ws.Range("C89").Value = 58
$ws.Range("D89").Value = 7240
$ws.Range("E89").Value = 840

# Row 92: Guinea -> Guinea
$ws.Range("B92").Value = 7664
$ws.Range("C92").Value = 89
$ws.Range("D92").Value = 6757
$ws.Range("E92").Value = 858

# Row 123: Mali -> Mali
$ws.Range("B123").Value = 2552
$ws.Range("C123").Value = 6
$ws.Range("D123").Value = 1954
$ws.Range("E123").Value = 474

# Row 137: Tunez -> Tunez
$ws.Range("B137").Value = 1642
$ws.Range("C137").Value = 41
$ws.Range("D137").Value = 1241
$ws.Range("E137").Value = 350

# Row 167: Burundi -> Comoras
$ws.Range("A167").Value = 'Comoras'
$ws.Range("B167").Value = 396
$ws.Range("C167").Value = 8
$ws.Range("D167").Value = 340
$ws.Range("E167").Value = 49
$ws.Range("H167").Value = 7

# Row 168: Comoras -> Burundi
$ws.Range("A168").Value = 'Burundi'
$ws.Range("B168").Value = 395
$ws.Range("D168").Value = 304
$ws.Range("E168").Value = 90
$ws.Range("H168").Value = 1

# Row 177: Camboya -> Aruba
$ws.Range("A177").Value = 'Aruba'
$ws.Range("B177").Value = 263
$ws.Range("C177").Value = 92
$ws.Range("D177").Value = 114
$ws.Range("E177").Value = 146
$ws.Range("H177").Value = 3

# Row 178: Trinidad yTobago -> Camboya
$ws.Range("A178").Value = 'Camboya'
$ws.Range("B178").Value = 243
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 210
$ws.Range("E178").Value = 33
$ws.Range("H178").Value = 0

# Row 179: Islas Caimanes -> Trinidad yTobago
$ws.Range("A179").Value = 'Trinidad yTobago'
$ws.Range("B179").Value = 207
$ws.Range("C179").Value = 8
$ws.Range("D179").Value = 135
$ws.Range("E179").Value = 64
$ws.Range("H179").Value = 8

# Row 180: Gibraltar -> Islas Caimanes
$ws.Range("A180").Value = 'Islas Caimanes'
$ws.Range("B180").Value = 203
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 202
$ws.Range("E180").Value = 0
$ws.Range("H180").Value = 1

# Row 181: Aruba -> Gibraltar
$ws.Range("A181").Value = 'Gibraltar'
$ws.Range("B181").Value = 190
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 184
$ws.Range("E181").Value = 6
$ws.Range("H181").Value = 0

# Row 202: Santa Lucia -> Timor Oriental
$ws.Range("A202").Value = 'Timor Oriental'

# Row 203: Timor Oriental -> Santa Lucia
$ws.Range("A203").Value = 'Santa Lucia'
